# "Remove form_id from basic forms"
# - settings sheet: delete the form_id column (B), shifting version/style/namespaces left
# - survey sheet: becomes non-active tab, settings sheet becomes active tab
# - survey sheet: simplify/merge the conditional-formatting ranges for columns A:G and C

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. settings sheet: remove the form_id column
# ---------------------------------------------------------------------------
$settings = $wb.Worksheets.Item("settings")

# Capture the existing header-row cell comments (by position) before we shift
# columns around, so we can re-attach them to the correct cells afterwards.
$commentFormTitle = $settings.Range("A1").Comment.Text()
$commentVersion = $settings.Range("C1").Comment.Text()
$commentStyle = $settings.Range("D1").Comment.Text()
$commentNamespaces = $settings.Range("E1").Comment.Text()

# Comments do not auto-shift with column deletion, so remove them all first.
$settings.Range("A1").Comment.Delete()
$settings.Range("B1").Comment.Delete()
$settings.Range("C1").Comment.Delete()
$settings.Range("D1").Comment.Delete()
$settings.Range("E1").Comment.Delete()

# Delete column B (form_id) entirely -- this shifts the version/style/namespaces
# columns (and their row-2 values) one column to the left automatically.
$settings.Columns.Item(2).Delete()

# Re-create the comments on their new cells.
$settings.Range("A1").AddComment($commentFormTitle)
$settings.Range("B1").AddComment($commentVersion)
$settings.Range("C1").AddComment($commentStyle)
$settings.Range("D1").AddComment($commentNamespaces)

# settings becomes the active sheet/tab, with A5 selected.
$settings.Activate()
$settings.Range("A5").Select()

# ---------------------------------------------------------------------------
# 2. survey sheet: simplify the conditional formatting ranges
# ---------------------------------------------------------------------------
$survey = $wb.Worksheets.Item("survey")

$wideFormat = $survey.Range("A2").FormatConditions
$wideFormat.Item(1).ModifyAppliesToRange($survey.Range("A2:G10000"))

$colCFormat = $survey.Range("C2").FormatConditions
$colCFormat.Item($colCFormat.Count).ModifyAppliesToRange($survey.Range("C2:C10000"))
